$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 5 new rows starting at row 2, shifting all existing data down
$insertRange = $ws.Range("A2:A6").EntireRow
$insertRange.Insert()

# New song entries to place into the newly inserted rows
$newSongs = @(
    @("Sunset Sons - VROL", "https://www.youtube.com/watch?v=vMIgQ36zhAg"),
    @("Passenger - Patient Love (Lyrics)", "https://www.youtube.com/watch?v=6LcKdxaSZVU"),
    @("Billie Eilish - lovely (Lyrics) ft. Khalid", "https://www.youtube.com/watch?v=8VLXHyHRXjc"),
    @("Arash feat.Helena - Angels Lullaby(Lyrics)", "https://www.youtube.com/watch?v=MeHCr0e-8vk"),
    @("Billie Eilish - Birds Of A Feather (Lyrics)", "https://www.youtube.com/watch?v=d5gf9dXbPi0")
)

$rowIndex = 2
foreach ($song in $newSongs) {
    $ws.Cells.Item($rowIndex, 2).Value = $song[0]
    $ws.Cells.Item($rowIndex, 3).Value = $song[1]
    $rowIndex++
}
